# Corrected write column error: the "Full_Name" values for rows 2-14 were
# mistakenly written into column I or J instead of column K. Move each
# value into column K and clear the original (wrong) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K is where the Full_Name value belongs for every data row.
$targetCol = 11  # K

for ($row = 2; $row -le 14; $row++) {
    # Find which cell among I (9) and J (10) holds the misplaced Full_Name value
    # (it is always whichever of the two is the last non-blank one in the row).
    $foundCol = $null
    foreach ($col in 9, 10) {
        $cell = $ws.Cells.Item($row, $col)
        if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") {
            $foundCol = $col
        }
    }

    if ($foundCol -ne $null) {
        $sourceCell = $ws.Cells.Item($row, $foundCol)
        $value = $sourceCell.Value2

        # Write value into the correct column K.
        $ws.Cells.Item($row, $targetCol).Value2 = $value

        # Clear the incorrectly placed cell.
        $sourceCell.ClearContents()
    }
}
